# Apply scheduled-runner data refresh updates to Hades_Profits workbook
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for a
# handful of rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets.

$wb = $excel.ActiveWorkbook

# Row 129 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 873.70734
$ws.Range("I129").Value = 500.4
$ws.Range("J129").Value = 994.129
$ws.Range("K129").Value = 1501.2
$ws.Range("L129").Value = 2982.387
$ws.Range("M129").Value = 3498.8
$ws.Range("N129").Value = -12982.387

# Row 135 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 34629
$ws.Range("I135").Value = 39582.848
$ws.Range("K135").Value = 356245.632
$ws.Range("M135").Value = -353710.632

# Row 137 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2942924.8
$ws.Range("I137").Value = 3449824
$ws.Range("J137").Value = 2908.6
$ws.Range("K137").Value = 10349472
$ws.Range("L137").Value = 8725.799999999999
$ws.Range("M137").Value = -10346922
$ws.Range("N137").Value = -13825.8

# Row 138 on ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2780527.8
$ws.Range("I138").Value = 1001.0526
$ws.Range("J138").Value = 4068601.2
$ws.Range("K138").Value = 3003.1578
$ws.Range("L138").Value = 12205803.6
$ws.Range("M138").Value = 2136.8422
$ws.Range("N138").Value = -12216083.6

# Row 61 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 23303262
$ws.Range("I61").Value = 30334018
$ws.Range("J61").Value = 101771.4
$ws.Range("K61").Value = 30334018
$ws.Range("L61").Value = 101771.4
$ws.Range("M61").Value = -30333806
$ws.Range("N61").Value = -102195.4

# Row 74 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 6633710
$ws.Range("I74").Value = 8966014
$ws.Range("J74").Value = 103260
$ws.Range("K74").Value = 8966014
$ws.Range("L74").Value = 103260
$ws.Range("M74").Value = -8965140
$ws.Range("N74").Value = -105008

# Row 77 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 6633710
$ws.Range("I77").Value = 8966014
$ws.Range("J77").Value = 103260
$ws.Range("K77").Value = 44830070
$ws.Range("L77").Value = 516300
$ws.Range("M77").Value = -44825702
$ws.Range("N77").Value = -525036

# Row 97 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 3677283
$ws.Range("I97").Value = 4464918
$ws.Range("K97").Value = 4464918
$ws.Range("M97").Value = -4464422

# Row 102 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 10205636
$ws.Range("I102").Value = 11906034
$ws.Range("J102").Value = 3249.5
$ws.Range("K102").Value = 11906034
$ws.Range("L102").Value = 3249.5
$ws.Range("M102").Value = -11904412
$ws.Range("N102").Value = -6493.5

# Row 122 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 5850117
$ws.Range("I122").Value = 2194.625
$ws.Range("K122").Value = 6583.875
$ws.Range("M122").Value = -4133.875

# Row 136 on ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 23303262
$ws.Range("I136").Value = 30334018
$ws.Range("J136").Value = 101771.4
$ws.Range("K136").Value = 91002054
$ws.Range("L136").Value = 305314.2
$ws.Range("M136").Value = -90999504
$ws.Range("N136").Value = -310414.2

# Row 20 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1136.0714
$ws.Range("I20").Value = 801.1667
$ws.Range("J20").Value = 1387.25
$ws.Range("K20").Value = 801.1667
$ws.Range("L20").Value = 1387.25
$ws.Range("M20").Value = -554.1667
$ws.Range("N20").Value = -1881.25

# Row 80 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 577.58826
$ws.Range("I80").Value = 220.25
$ws.Range("J80").Value = 895.2222
$ws.Range("K80").Value = 220.25
$ws.Range("L80").Value = 895.2222
$ws.Range("M80").Value = 777.75
$ws.Range("N80").Value = -2891.2222

# Row 83 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 577.58826
$ws.Range("I83").Value = 220.25
$ws.Range("J83").Value = 895.2222
$ws.Range("K83").Value = 1101.25
$ws.Range("L83").Value = 4476.111
$ws.Range("M83").Value = 3890.75
$ws.Range("N83").Value = -14460.111

# Row 134 on BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2298.9333
$ws.Range("I134").Value = 1786.7872
$ws.Range("K134").Value = 5360.3616
$ws.Range("M134").Value = -2825.3616

# Row 7 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 356.2
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 420.25
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 420.25
$ws.Range("M7").Value = 13
$ws.Range("N7").Value = -646.25

# Row 31 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2213.5293
$ws.Range("I31").Value = 1217.9744
$ws.Range("J31").Value = 5449.0835
$ws.Range("K31").Value = 1217.9744
$ws.Range("L31").Value = 5449.0835
$ws.Range("M31").Value = -922.9744000000001
$ws.Range("N31").Value = -6039.0835

# Row 34 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2213.5293
$ws.Range("I34").Value = 1217.9744
$ws.Range("J34").Value = 5449.0835
$ws.Range("K34").Value = 1217.9744
$ws.Range("L34").Value = 5449.0835
$ws.Range("M34").Value = -1015.9744
$ws.Range("N34").Value = -5853.0835

# Row 58 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 17545186
$ws.Range("I58").Value = 20834488
$ws.Range("J58").Value = 2244.4443
$ws.Range("K58").Value = 20834488
$ws.Range("L58").Value = 2244.4443
$ws.Range("M58").Value = -20834285
$ws.Range("N58").Value = -2650.4443

# Row 107 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 428.83334
$ws.Range("I107").Value = 405.75
$ws.Range("J107").Value = 475
$ws.Range("K107").Value = 405.75
$ws.Range("L107").Value = 475
$ws.Range("M107").Value = 1514.25
$ws.Range("N107").Value = -4315

# Row 132 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 19655.564
$ws.Range("I132").Value = 1184.2391
$ws.Range("J132").Value = 114064.555
$ws.Range("K132").Value = 3552.7173
$ws.Range("L132").Value = 342193.665
$ws.Range("M132").Value = -1022.7173
$ws.Range("N132").Value = -347253.665

# Row 134 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 24861.164
$ws.Range("I134").Value = 1877.359
$ws.Range("J134").Value = 114498
$ws.Range("K134").Value = 5632.076999999999
$ws.Range("L134").Value = 343494
$ws.Range("M134").Value = -3097.076999999999
$ws.Range("N134").Value = -348564

# Row 136 on CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 17545186
$ws.Range("I136").Value = 20834488
$ws.Range("J136").Value = 2244.4443
$ws.Range("K136").Value = 62503464
$ws.Range("L136").Value = 6733.3329
$ws.Range("M136").Value = -62500914
$ws.Range("N136").Value = -11833.3329

# Row 14 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 835.7037
$ws.Range("I14").Value = 835.7037
$ws.Range("K14").Value = 2507.1111
$ws.Range("M14").Value = -2334.1111

# Row 131 on CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1059.5333
$ws.Range("I131").Value = 793.3333
$ws.Range("J131").Value = 1073.5438
$ws.Range("K131").Value = 2379.9999
$ws.Range("L131").Value = 3220.6314
$ws.Range("M131").Value = 2660.0001
$ws.Range("N131").Value = -13300.6314

# Row 100 on GSM: H,J,L -> 0; clear N (cell removed)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

# Row 123 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 26801.455
$ws.Range("J123").Value = 26801.455
$ws.Range("L123").Value = 26801.455
$ws.Range("N123").Value = -31701.455

# Row 132 on GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 39861.598
$ws.Range("I132").Value = 28223.459
$ws.Range("J132").Value = 68569
$ws.Range("K132").Value = 84670.37699999999
$ws.Range("L132").Value = 205707
$ws.Range("M132").Value = -82140.37699999999
$ws.Range("N132").Value = -210767

# Row 132 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 39985.4
$ws.Range("I132").Value = 16953.47
$ws.Range("J132").Value = 170499.67
$ws.Range("K132").Value = 50860.41
$ws.Range("L132").Value = 511499.01
$ws.Range("M132").Value = -48330.41
$ws.Range("N132").Value = -516559.01

# Row 136 on LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 49934.69
$ws.Range("I136").Value = 27462.025
$ws.Range("J136").Value = 263425
$ws.Range("K136").Value = 82386.07500000001
$ws.Range("L136").Value = 790275
$ws.Range("M136").Value = -79836.07500000001
$ws.Range("N136").Value = -795375

# Row 97 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H97").Value = 53000
$ws.Range("J97").Value = 53000
$ws.Range("L97").Value = 53000
$ws.Range("N97").Value = -54982

# Row 132 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 62209.938
$ws.Range("I132").Value = 48720.715
$ws.Range("K132").Value = 146162.145
$ws.Range("M132").Value = -143632.145

# Row 136 on WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 33052.793
$ws.Range("I136").Value = 19830.584
$ws.Range("J136").Value = 103130.5
$ws.Range("K136").Value = 59491.75199999999
$ws.Range("L136").Value = 309391.5
$ws.Range("M136").Value = -56941.75199999999
$ws.Range("N136").Value = -314491.5

